$wb = $excel.ActiveWorkbook

# Sheet 1: 展览
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(2, 6).Value = 93
$ws.Cells.Item(3, 6).Value = 123
$ws.Cells.Item(5, 6).Value = 362
$ws.Cells.Item(6, 6).Value = 547
$ws.Cells.Item(7, 6).Value = 1535
$ws.Cells.Item(8, 6).Value = 11670
$ws.Cells.Item(9, 6).Value = 196
$ws.Cells.Item(10, 6).Value = 76
$ws.Cells.Item(11, 6).Value = 122
$ws.Cells.Item(12, 6).Value = 2096
$ws.Cells.Item(15, 6).Value = 7
$ws.Cells.Item(16, 6).Value = 235
$ws.Cells.Item(18, 6).Value = 1190
$ws.Cells.Item(19, 6).Value = 167
$ws.Cells.Item(20, 6).Value = 247
$ws.Cells.Item(21, 6).Value = 736
$ws.Cells.Item(23, 6).Value = 264
$ws.Cells.Item(24, 6).Value = 2402
$ws.Cells.Item(25, 6).Value = 724
$ws.Cells.Item(26, 6).Value = 3543
$ws.Cells.Item(27, 6).Value = 1057
$ws.Cells.Item(28, 6).Value = 792
$ws.Cells.Item(32, 6).Value = 978
$ws.Cells.Item(33, 6).Value = 35
$ws.Cells.Item(34, 6).Value = 61
$ws.Cells.Item(36, 6).Value = 18
$ws.Cells.Item(38, 6).Value = 9
$ws.Cells.Item(39, 6).Value = 2546
$ws.Cells.Item(40, 6).Value = 4394
$ws.Cells.Item(41, 6).Value = 5458
$ws.Cells.Item(43, 6).Value = 113
$ws.Cells.Item(44, 6).Value = 152
$ws.Cells.Item(45, 6).Value = 262
$ws.Cells.Item(46, 6).Value = 63
$ws.Cells.Item(47, 6).Value = 28
$ws.Cells.Item(48, 6).Value = 4091
$ws.Cells.Item(49, 6).Value = 98

# Sheet 2: 演出
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(3, 6).Value = 4149
$ws.Cells.Item(11, 6).Value = 646
$ws.Cells.Item(15, 6).Value = 3
$ws.Cells.Item(16, 6).Value = 9

# Sheet 3: 本地生活
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(3, 6).Value = 421

# Sheet 4: 全部类型
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(3, 6).Value = 421
$ws.Cells.Item(5, 6).Value = 93
$ws.Cells.Item(6, 6).Value = 123
$ws.Cells.Item(8, 6).Value = 362
$ws.Cells.Item(9, 6).Value = 547
$ws.Cells.Item(10, 6).Value = 11670
$ws.Cells.Item(11, 6).Value = 76
$ws.Cells.Item(12, 6).Value = 122
$ws.Cells.Item(13, 6).Value = 2096
$ws.Cells.Item(16, 6).Value = 7
$ws.Cells.Item(18, 6).Value = 1190
$ws.Cells.Item(19, 6).Value = 167
$ws.Cells.Item(20, 6).Value = 247
$ws.Cells.Item(21, 6).Value = 4149
$ws.Cells.Item(23, 6).Value = 736
$ws.Cells.Item(24, 6).Value = 724
$ws.Cells.Item(26, 6).Value = 792
$ws.Cells.Item(31, 6).Value = 978
$ws.Cells.Item(32, 6).Value = 61
$ws.Cells.Item(34, 6).Value = 18
$ws.Cells.Item(36, 6).Value = 4394
$ws.Cells.Item(39, 6).Value = 113
$ws.Cells.Item(40, 6).Value = 152
$ws.Cells.Item(41, 6).Value = 262
$ws.Cells.Item(43, 6).Value = 3
$ws.Cells.Item(44, 6).Value = 63
$ws.Cells.Item(45, 6).Value = 4091
$ws.Cells.Item(48, 6).Value = 98
